# Project Sample Project is saved. Author: admin. Type: SAVE.
#
# Cell B11 on the "Rules" sheet currently holds the text "R40" (row for the
# 4th greeting-time bracket). The saved edit replaces that text with the
# literal text "1", keeping the cell's existing formatting/style untouched
# and still stored as a (shared) text string, not a number.
#
# A plain  Range.Value = "1"  assignment would be auto-coerced to the
# number 1 by Excel's normal type inference, and prefixing with an
# apostrophe (forcing text) stamps the cell with a "quote prefix" which
# bumps it onto a brand new cell style. Neither preserves the original
# file shape, so instead we go through a formula that evaluates to the
# text "1", then collapse that formula down to its literal value with a
# values-only paste - this keeps the result typed as text while leaving
# the cell's style/formatting completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")
$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues - drop the formula, keep the text value + original style
